# Insert the new "roteiro" content (requisitos, modelagem, especificacao de
# componentes, etc.) right after the existing "Quadro comparativo entre as
# ferramentas de ETL - revisao de literatura." paragraph and the blank
# paragraph that follows it, while leaving the trailing bookmark paragraph
# ("_GoBack") untouched.

$d = $word.ActiveDocument

# Locate the anchor paragraph by its text (robust to paragraph-index drift).
$findRange = $d.Content
$found = $findRange.Find.Execute("Quadro comparativo entre as ferramentas de ETL", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Anchor paragraph 'Quadro comparativo...' not found"
}

# $findRange now covers just the matched text; expand to its paragraph, then
# to the following (blank) paragraph, and insert right after that one.
$anchorPara = $findRange.Paragraphs(1)
$afterAnchor = $d.Range($anchorPara.Range.End, $anchorPara.Range.End)
$blankPara = $afterAnchor.Paragraphs(1)
$insertPoint = $d.Range($blankPara.Range.End, $blankPara.Range.End)

$newContentXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Para o nosso trabalho propomos elencar os requisitos de software do ETL4NoSQL, as modelagens do domínio e especificação, a especificação dos componentes, interfaces de programação, o estudo experimental de software e as especializações a partir do ETL4NoSQL. Estarei explicando esses itens com mais detalhe a seguir.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">Os requisitos de software elencados para o ETL4NoSQL foram: </w:t></w:r></w:p><w:p/><w:p><w:r><w:t>A entidade Modelagem faz a leitura de uma o</w:t></w:r><w:r><w:t xml:space="preserve">u mais Fontes de dados, executa </w:t></w:r><w:r><w:t>nenhuma ou muitas Operações. As Operações por sua</w:t></w:r><w:r><w:t xml:space="preserve"> vez são processadas por um ou </w:t></w:r><w:r><w:t>mais Processamentos de forma a serem P</w:t></w:r><w:r><w:t>rocessamentos Centralizados ou processamentos</w:t></w:r><w:r><w:t xml:space="preserve"> Distribuídos. E por fim, a Mo</w:t></w:r><w:r><w:t xml:space="preserve">delagem carrega o resultado das </w:t></w:r><w:r><w:t>Operações processadas à um ou mais Destinos</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Os casos de uso expressam as funcionalidades fundamentais ao desenvolvimento e uso do ETL4NoSQL. Eles permitem ao programador a visão do que é imprescindível ao implementar e determinar as interfaces de sistema e operações do framework.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Um conjunto de casos de uso foram identificados tais como: Ler fonte de dados, </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Escrever</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> no destino, Modelar dados, Executar </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>operação</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> e Processar </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>operações</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. </w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">Ao construir o modelo comportamental é </w:t></w:r><w:r><w:t>possível</w:t></w:r><w:r><w:t xml:space="preserve"> id</w:t></w:r><w:r><w:t>entificar os conceitos com com</w:t></w:r><w:r><w:t xml:space="preserve">portamentos mais relevantes para o </w:t></w:r><w:r><w:t>negócio</w:t></w:r><w:r><w:t xml:space="preserve">, bem como os estados e eventos que disparam as transições entre os estados (SOUZA GIMENES; HUZITA (2005)). Dessa forma, </w:t></w:r><w:r><w:t xml:space="preserve">apresentamos </w:t></w:r><w:r><w:t>o diagrama de estados do ETL4NoSQL, nele podemos ver as transições de leitura da fonte de dados, validação e identificação dos dados, assim como o tratamento caso os dados não possam ser identificados. Subsequente a isso, podemos ver as transições do armazenamento dos dados para o processamento, a criação dos processos de ETL,</w:t></w:r><w:r><w:t xml:space="preserve"> a escolha da forma de process</w:t></w:r><w:r><w:t>amento, execução das operações, e também o tratamento para as operações que não puderem ser executadas. Finalmente, a transição da carga dos dados pode ser feita na base de destino seguido da mensagem de tratamento, caso haja sucesso ou não na execução.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">A modelagem de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>especificação</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> visa definir, em um </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>nível</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> alto de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>abstração</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, os </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>serviços</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> oferecidos pelos componentes. Dessa forma, é </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>possível</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> determinar a arquitetura e especificar os seus componentes. É importante dar </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ênfase</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>especificação</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> das interfaces, pois isso contribui para uma clara </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>separação</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> entre os componentes, e </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>também</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, para assegurar o </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>princípio</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> de encapsul</w:t></w:r><w:r><w:t xml:space="preserve">amento de dados e comportamento. </w:t></w:r><w:r><w:t xml:space="preserve">CHEESMAN; DANIELS (2001) divide a modelagem da </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>especificação</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">em </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>três</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>estágios</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>: Identificação de componentes, Interfaces de Sistemas e Interfaces de Negócio.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">As interfaces de sistemas são identificadas por meio dos casos de uso, mencionados na fase de licitação de requisitos. Assim, identificamos as interfaces: </w:t></w:r></w:p><w:p/><w:p><w:r><w:t>As interfaces de negócio são identificadas por meio do modelo conceitual, cada entidade identificada no modelo conceitual é uma interface de negócio.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">Assim, após definir as interfaces de sistemas e de negócio podemos </w:t></w:r><w:r><w:t>determinar a arquitetura dos componentes.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">Depois de </w:t></w:r><w:r><w:t xml:space="preserve">definir as operações foi possível criar o diagrama de especificação das interfaces. </w:t></w:r><w:r><w:t>Essas especificações referem-se ao uso dos componentes. Elas representam o que o usuário precisa saber sobre os componentes.</w:t></w:r><w:r><w:t xml:space="preserve"> Aqui nós temos as interfaces de negócio identificadas, as interações entre os componentes e as operações.</w:t></w:r></w:p><w:p/>'

$null = $insertPoint.InsertXML($newContentXml)
